$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at the top of the "Ajo" data block (rows 528-529),
# pushing the existing rows 528-566 down to 530-568.
$ws.Rows("528:529").Insert()

# ---- New row 528 ----
$ws.Range("A528").Value = 4
$ws.Range("B528").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C528").Value = "Los Lagos"
$ws.Range("D528").Value = 45265
$ws.Range("E528").Value = 10
$ws.Range("F528").Value = 100112003
$ws.Range("G528").Value = "Ajo"
$ws.Range("H528").Value = "Chino"
$ws.Range("I528").Value = "Primera"
$ws.Range("J528").Value = 250
$ws.Range("K528").Value = 26000
$ws.Range("L528").Value = 26000
$ws.Range("M528").Value = 26000
$ws.Range("N528").Value = "`$/caja 10 kilos"
$ws.Range("O528").Value = "China"
$ws.Range("P528").Value = 2600
$ws.Range("Q528").Value = 10
$ws.Range("R528").Value = "Hortaliza"

# ---- New row 529 ----
$ws.Range("A529").Value = 4
$ws.Range("B529").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C529").Value = "Los Lagos"
$ws.Range("D529").Value = 45265
$ws.Range("E529").Value = 10
$ws.Range("F529").Value = 100112003
$ws.Range("G529").Value = "Ajo"
$ws.Range("H529").Value = "Chino"
$ws.Range("I529").Value = "Primera"
$ws.Range("J529").Value = 100
$ws.Range("K529").Value = 27000
$ws.Range("L529").Value = 27000
$ws.Range("M529").Value = 27000
$ws.Range("N529").Value = "`$/malla 10 kilos"
$ws.Range("O529").Value = "China"
$ws.Range("P529").Value = 2700
$ws.Range("Q529").Value = 10
$ws.Range("R529").Value = "Hortaliza"
